# Restructure runs with spell-check proofErr markers (w:proofErr) around
# English/technical terms embedded in the Italian text, per the source edit,
# and append the new trailing todo item about the Friends/pending page.
$d = $word.ActiveDocument

$d.Paragraphs(1).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Citazioneintensa"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>ToDo</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>things</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>') | Out-Null
$d.Paragraphs(2).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Quando cambio utente, la funzione </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>useEffect</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> che crea il </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>context</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> nel Drawer.js stampa due volte il valore di console.log(), prova a mettere una </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>useCallback</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>()</w:t></w:r><w:r><w:t xml:space="preserve">   </w:t></w:r></w:p>') | Out-Null
$d.Paragraphs(3).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>STYLE_</w:t></w:r><w:r><w:t>Utilizzare</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> un’applicazione per creare lo stile delle pagine alla fine delle implementazioni funzionali.</w:t></w:r></w:p>') | Out-Null
$d.Paragraphs(4).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">STYLE_ la </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>SearchBar</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> della pagina amici meglio nasconderla sotto l’</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>header</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> e mostrarla a scorrimento dall’alto verso il basso quando viene premuto un bottone di ricerca manualmente</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p>') | Out-Null
$d.Paragraphs(5).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Aggiornare il </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Modal</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> all’interno del component </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>MultiSelection</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> in modo tale da avere uno scope più generale. Modificarlo permettendo all’utente di passare dei &lt;</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Toucahble</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>/&gt; come figli in questo modo :</w:t></w:r><w:r><w:br/></w:r><w:r><w:br/><w:t>&lt;</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>MultiSelection</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>&gt;</w:t></w:r></w:p>') | Out-Null
$d.Paragraphs(6).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:ind w:left="2124"/></w:pPr><w:r><w:br/><w:t>&lt;</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Touhable</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>&gt;</w:t></w:r><w:r><w:br/><w:t>…</w:t></w:r><w:r><w:br/><w:t>&lt;/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Touchable</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>&gt;</w:t></w:r><w:r><w:br/></w:r></w:p>') | Out-Null
$d.Paragraphs(7).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:left="708" w:firstLine="708"/></w:pPr><w:r><w:t>&lt;/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Multiselection</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>&gt;</w:t></w:r></w:p>') | Out-Null
$d.Paragraphs(9).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">L’API che riempie gli amici nella </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Flatlist</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> mostra solo i primi 50 amici, considera di estendere gli amici ogni qualvolta l’utente raggiunge il fondo della </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Flatlist</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> di 10 persone aggiuntive. Lo puoi implementare passando un parametro (un numero intero) all’ API </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>getPossibleFriendsBySimilarUsername</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> per dirgli che deve iniziare da quel valore a cercare nuovi amici</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p>') | Out-Null
$d.Paragraphs(10).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Gestire Pagina di </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>resgistrazione</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p>') | Out-Null
$d.Paragraphs(11).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Ogni volta che viene creato un nuovo utente nel database, far partire un trigger che crea un analogo documento nella </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>collection</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>notifications</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”.</w:t></w:r></w:p>') | Out-Null

# Append the brand-new final todo paragraph (doesn't exist in the source doc yet)
$lastIndex = $d.Paragraphs.Count
$d.Paragraphs($lastIndex).Range.InsertParagraphAfter() | Out-Null
$d.Paragraphs($d.Paragraphs.Count).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Nella sezione Friends </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>inlcudere</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> una pagina di amici in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pending</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p>') | Out-Null

